# TestData.xlsx update
# - hrms_id test value changes from 98 -> 75 on every sheet that carries it
# - OrderHeaderPage / OrderDetailPage occasion date moves from 15-Nov-2022 /
#   30-Nov-2022 to reflect a new occasionmonth "Dec" (and, on OrderDetailPage,
#   occasiondate "10")
# - Active cell / selection on a few sheets is updated to match where the
#   user left off editing

$wb = $excel.ActiveWorkbook

$wsSalesLogin   = $wb.Worksheets.Item("SalesLogin")
$wsAddNew       = $wb.Worksheets.Item("AddNewCustomerDetails")
$wsSearch       = $wb.Worksheets.Item("SearchCustomerDetails")
$wsOrderHeader  = $wb.Worksheets.Item("OrderHeaderPage")
$wsOrderDetail  = $wb.Worksheets.Item("OrderDetailPage")

# SalesLogin: hrms_id 98 -> 75
$wsSalesLogin.Range("A2").Value = "75"

# AddNewCustomerDetails: hrms_id 98 -> 75, active cell moves to B2
$wsAddNew.Range("A2").Value = "75"
$wsAddNew.Range("B2").Select()

# SearchCustomerDetails: hrms_id 98 -> 75 on every data row
$wsSearch.Range("A2").Value = "75"
$wsSearch.Range("A3").Value = "75"
$wsSearch.Range("A4").Value = "75"
$wsSearch.Range("A5").Value = "75"
$wsSearch.Range("A6").Value = "75"

# OrderHeaderPage: hrms_id 98 -> 75, occasionmonth Nov -> Dec, active cell -> H2
$wsOrderHeader.Range("A2").Value = "75"
$wsOrderHeader.Range("H2").Value = "Dec"
$wsOrderHeader.Range("H2").Select()

# OrderDetailPage: hrms_id 98 -> 75, occasionmonth Nov -> Dec,
# occasiondate 30 -> 10, active cell -> J2
$wsOrderDetail.Range("A2").Value = "75"
$wsOrderDetail.Range("H2").Value = "Dec"
$wsOrderDetail.Range("J2").Value = "10"
$wsOrderDetail.Range("J2").Select()
